$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 122, pushing existing rows 122:225 down to 123:226
$ws.Rows.Item(122).Insert()

# Populate the newly inserted row 122 with the new weekly record
$ws.Cells.Item(122, 1).Value2 = 5
$ws.Cells.Item(122, 2).Value2 = "Macroferia Regional de Talca"
$ws.Cells.Item(122, 3).Value2 = "Maule"
$ws.Cells.Item(122, 4).Value2 = 44651
$ws.Cells.Item(122, 5).Value2 = 7
$ws.Cells.Item(122, 6).Value2 = 100112024
$ws.Cells.Item(122, 7).Value2 = "Choclo"
$ws.Cells.Item(122, 8).Value2 = "Choclero"
$ws.Cells.Item(122, 9).Value2 = "Primera"
$ws.Cells.Item(122, 10).Value2 = 20000
$ws.Cells.Item(122, 11).Value2 = 200
$ws.Cells.Item(122, 12).Value2 = 200
$ws.Cells.Item(122, 13).Value2 = 200
$ws.Cells.Item(122, 14).Value2 = "`$/unidad"
$ws.Cells.Item(122, 15).Value2 = "Región del Maule"
$ws.Cells.Item(122, 16).Value2 = 200
$ws.Cells.Item(122, 17).Value2 = 1
$ws.Cells.Item(122, 18).Value2 = "Hortaliza"
